# The commit adds one new weekly price record for
# "Vega Monumental Concepción - Ajo" that belongs right after the existing
# row for 2021-02-11 (serial 44238), i.e. at sheet row 86 (row 2 is the
# first data row). Inserting a whole row there shifts every subsequent
# record down by one (old row 86 -> new row 87, ..., old row 183 -> new
# row 184), which matches the observed diff (dimension grows from R183 to
# R184, and each D/J/K/L/M/P value reappears one row lower than before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 86, pushing rows 86:183 down to 87:184
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A86").Value2 = 11
$ws.Range("B86").Value2 = "Vega Monumental Concepción"
$ws.Range("C86").Value2 = "Bíobío"
$ws.Range("D86").Value2 = 44740
$ws.Range("E86").Value2 = 8
$ws.Range("F86").Value2 = 100112003
$ws.Range("G86").Value2 = "Ajo"
$ws.Range("H86").Value2 = "Chino"
$ws.Range("I86").Value2 = "Primera"
$ws.Range("J86").Value2 = 270
$ws.Range("K86").Value2 = 15000
$ws.Range("L86").Value2 = 16000
$ws.Range("M86").Value2 = 15444
$ws.Range("N86").Value2 = "$/caja 10 kilos"
$ws.Range("O86").Value2 = "China"
$ws.Range("P86").Value2 = 1544
$ws.Range("Q86").Value2 = 10
$ws.Range("R86").Value2 = "Hortaliza"
